{"js": "const replacements = [\n  [\"2025-11-24 Monday\", \"2025-11-25 Tuesday\"],\n  [\"14\u00d711=\", \"60\u00d748=\"],\n  [\"95\u00d798=\", \"35\u00d752=\"],\n  [\"88\u00d744=\", \"35\u00d724=\"],\n  [\"96\u00d767=\", \"58\u00d793=\"],\n  [\"90\u00d738=\", \"60\u00d745=\"],\n  [\"95\u00d796=\", \"48\u00d716=\"],\n  [\"17\u00d721=\", \"97\u00d778=\"],\n  [\"53\u00d738=\", \"53\u00d762=\"],\n  [\"59\u00d772=\", \"55\u00d776=\"],\n  [\"99\u00d772=\", \"96\u00d761=\"],\n  [\"27\u00d730=\", \"79\u00d755=\"],\n  [\"42\u00d772=\", \"89\u00d712=\"],\n  [\"84\u00d712=\", \"40\u00d736=\"],\n  [\"79\u00d724=\", \"36\u00d736=\"],\n  [\"27\u00d738=\", \"75\u00d721=\"],\n  [\"92\u00d787=\", \"36\u00d735=\"],\n  [\"54\u00d716=\", \"72\u00d717=\"],\n  [\"48\u00d759=\", \"75\u00d791=\"],\n  [\"41\u00d793=\", \"23\u00d788=\"],\n  [\"61\u00d715=\", \"55\u00d763=\"],\n  [\"89\u00d787=\", \"11\u00d713=\"],\n  [\"67\u00d760=\", \"17\u00d754=\"],\n  [\"14\u00d780=\", \"83\u00d752=\"],\n  [\"12\u00d790=\", \"24\u00d790=\"],\n  [\"40\u00d733=\", \"68\u00d790=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $findText, $replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text $d \"2025-11-24 Monday\" \"2025-11-25 Tuesday\"\nReplace-Text $d \"14\u00d711=\" \"60\u00d748=\"\nReplace-Text $d \"95\u00d798=\" \"35\u00d752=\"\nReplace-Text $d \"88\u00d744=\" \"35\u00d724=\"\nReplace-Text $d \"96\u00d767=\" \"58\u00d793=\"\nReplace-Text $d \"90\u00d738=\" \"60\u00d745=\"\nReplace-Text $d \"95\u00d796=\" \"48\u00d716=\"\nReplace-Text $d \"17\u00d721=\" \"97\u00d778=\"\nReplace-Text $d \"53\u00d738=\" \"53\u00d762=\"\nReplace-Text $d \"59\u00d772=\" \"55\u00d776=\"\nReplace-Text $d \"99\u00d772=\" \"96\u00d761=\"\nReplace-Text $d \"27\u00d730=\" \"79\u00d755=\"\nReplace-Text $d \"42\u00d772=\" \"89\u00d712=\"\nReplace-Text $d \"84\u00d712=\" \"40\u00d736=\"\nReplace-Text $d \"79\u00d724=\" \"36\u00d736=\"\nReplace-Text $d \"27\u00d738=\" \"75\u00d721=\"\nReplace-Text $d \"92\u00d787=\" \"36\u00d735=\"\nReplace-Text $d \"54\u00d716=\" \"72\u00d717=\"\nReplace-Text $d \"48\u00d759=\" \"75\u00d791=\"\nReplace-Text $d \"41\u00d793=\" \"23\u00d788=\"\nReplace-Text $d \"61\u00d715=\" \"55\u00d763=\"\nReplace-Text $d \"89\u00d787=\" \"11\u00d713=\"\nReplace-Text $d \"67\u00d760=\" \"17\u00d754=\"\nReplace-Text $d \"14\u00d780=\" \"83\u00d752=\"\nReplace-Text $d \"12\u00d790=\" \"24\u00d790=\"\nReplace-Text $d \"40\u00d733=\" \"68\u00d790=\"\n"}
